# Applies the "add EL_Price and modify forecast" edit to the CTS_user_set workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Rework the SPEC_EN_EMPLOYEE "ELEC / H2 / HEAT" forecast rows
#    for every one of the 7 technology groups (rows 10-18, 19-27,
#    28-36, 37-45, 46-54, 55-63, 64-72): the Unit switches from
#    "take from hist!!" to "GWh/tsd. employees" and a real formula
#    (=1/1000000) replaces the placeholder Factor text.
# ------------------------------------------------------------------
$groupStarts = @(10, 19, 28, 37, 46, 55, 64)

foreach ($g in $groupStarts) {
    $elecRow = $g
    $h2Row = $g + 1
    $heatRow = $g + 2

    # --- ELEC row: plain formula ---
    $ws.Range("H$elecRow").Value = "GWh/tsd. employees"
    $ws.Range("H$elecRow").Style = "Normal"
    $ws.Range("I$elecRow").Formula = "=1/1000000"

    # --- H2 + HEAT rows: shared formula over the two-cell range ---
    $ws.Range("H$h2Row").Value = "GWh/tsd. employees"
    $ws.Range("H$h2Row").Style = "Normal"

    if ($heatRow -eq 66) {
        $ws.Range("H$heatRow").Value = "GWh/tsd. Employees"
    } else {
        $ws.Range("H$heatRow").Value = "GWh/tsd. employees"
    }
    $ws.Range("H$heatRow").Style = "Normal"

    $ws.Range("I" + $h2Row + ":I" + $heatRow).Formula = "=1/1000000"
}

# The very first group also picked up an (empty) styled N10 cell.
$ws.Range("N10").Style = "Standard 3"

# ------------------------------------------------------------------
# 2) New column AC ("Lower limit") with 0 for every data row.
# ------------------------------------------------------------------
$ws.Range("AC1").Value = "Lower limit"
for ($r = 2; $r -le 72; $r++) {
    $ws.Cells.Item($r, 29).Value = 0
}

# ------------------------------------------------------------------
# 3) Turn the AutoFilter off (keeps the hidden _FilterDatabase name,
#    but drops the <autoFilter> element) and grow the hidden
#    _FilterDatabase defined name to cover the full data range.
# ------------------------------------------------------------------
$filterOn = $ws.AutoFilterMode
if ($filterOn) {
    $ws.AutoFilterMode = $false
}
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Data!_FilterDatabase") {
        $n.RefersTo = "=Data!`$A`$1:`$AU`$72"
    }
}

# ------------------------------------------------------------------
# 4) Selection moves to N14.
# ------------------------------------------------------------------

[void]$ws.Range("N14").Select()

Write-Host "Edit applied"
